$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Fri Dec  8 10_41_10 2023"
$ws.Range("B4").Value = "loclexyz99"
$ws.Range("C4").Value = -3

$ws.Range("A5").Value = "Fri Dec  8 10_42_43 2023"
$ws.Range("B5").Value = "loclexyz99"
$ws.Range("C5").Value = -1
